{"js": "// Update the date line and each division-problem cell in the worksheet\n// table to the values from the new day's generated output.\nconst replacements = [\n  [\"2024-09-26 Thursday\", \"2024-09-27 Friday\"],\n  [\"356\u00f72=\", \"643\u00f72=\"],\n  [\"783\u00f72=\", \"347\u00f74=\"],\n  [\"244\u00f72=\", \"222\u00f75=\"],\n  [\"178\u00f76=\", \"789\u00f78=\"],\n  [\"170\u00f72=\", \"539\u00f73=\"],\n  [\"190\u00f72=\", \"482\u00f72=\"],\n  [\"945\u00f72=\", \"305\u00f75=\"],\n  [\"195\u00f74=\", \"354\u00f73=\"],\n  [\"104\u00f73=\", \"235\u00f72=\"],\n  [\"334\u00f79=\", \"300\u00f72=\"],\n  [\"783\u00f74=\", \"286\u00f75=\"],\n  [\"579\u00f75=\", \"978\u00f74=\"],\n  [\"718\u00f79=\", \"269\u00f78=\"],\n  [\"221\u00f73=\", \"707\u00f72=\"],\n  [\"335\u00f79=\", \"329\u00f78=\"],\n  [\"139\u00f79=\", \"163\u00f79=\"],\n  [\"520\u00f77=\", \"783\u00f76=\"],\n  [\"874\u00f74=\", \"966\u00f79=\"],\n  [\"699\u00f73=\", \"689\u00f73=\"],\n  [\"918\u00f78=\", \"860\u00f75=\"],\n  [\"689\u00f75=\", \"650\u00f76=\"],\n  [\"243\u00f75=\", \"747\u00f78=\"],\n  [\"628\u00f75=\", \"928\u00f78=\"],\n  [\"757\u00f77=\", \"490\u00f75=\"],\n  [\"337\u00f77=\", \"471\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each division-problem cell in the worksheet\n# table to the values from the new day's generated output.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-26 Thursday\", \"2024-09-27 Friday\"),\n    @(\"356\u00f72=\", \"643\u00f72=\"),\n    @(\"783\u00f72=\", \"347\u00f74=\"),\n    @(\"244\u00f72=\", \"222\u00f75=\"),\n    @(\"178\u00f76=\", \"789\u00f78=\"),\n    @(\"170\u00f72=\", \"539\u00f73=\"),\n    @(\"190\u00f72=\", \"482\u00f72=\"),\n    @(\"945\u00f72=\", \"305\u00f75=\"),\n    @(\"195\u00f74=\", \"354\u00f73=\"),\n    @(\"104\u00f73=\", \"235\u00f72=\"),\n    @(\"334\u00f79=\", \"300\u00f72=\"),\n    @(\"783\u00f74=\", \"286\u00f75=\"),\n    @(\"579\u00f75=\", \"978\u00f74=\"),\n    @(\"718\u00f79=\", \"269\u00f78=\"),\n    @(\"221\u00f73=\", \"707\u00f72=\"),\n    @(\"335\u00f79=\", \"329\u00f78=\"),\n    @(\"139\u00f79=\", \"163\u00f79=\"),\n    @(\"520\u00f77=\", \"783\u00f76=\"),\n    @(\"874\u00f74=\", \"966\u00f79=\"),\n    @(\"699\u00f73=\", \"689\u00f73=\"),\n    @(\"918\u00f78=\", \"860\u00f75=\"),\n    @(\"689\u00f75=\", \"650\u00f76=\"),\n    @(\"243\u00f75=\", \"747\u00f78=\"),\n    @(\"628\u00f75=\", \"928\u00f78=\"),\n    @(\"757\u00f77=\", \"490\u00f75=\"),\n    @(\"337\u00f77=\", \"471\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nWrite-Output \"done\"\n"}
